# Apply cell value updates for the crypto price table refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.078.29"
$ws.Range("E2").Value = "  -3.18%  "

$ws.Range("D3").Value = "2.559.29"
$ws.Range("E3").Value = "  -3.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.25%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("E13").Value = "  +7.91%  "

$ws.Range("D14").Value = "2.949.53"
$ws.Range("E14").Value = "  -3.54%  "

$ws.Range("D15").Value = "2.560.10"
$ws.Range("E15").Value = "  -3.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.884"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.98%  "

$ws.Range("D18").Value = "43.058.25"
$ws.Range("E18").Value = "  -3.33%  "

$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.50%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  -0.94%  "

$ws.Range("E21").Value = "  -2.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.05%  "

$ws.Range("E24").Value = "  -2.02%  "

$ws.Range("E25").Value = "  -6.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.67%  "

$ws.Range("E35").Value = "  -6.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0798"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.67%  "

$ws.Range("E39").Value = "  -2.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +38.85%  "

$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").Value = "2.107.06"
$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.69%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.805.71"
$ws.Range("E50").Value = "  -3.56%  "

$ws.Range("E51").Value = "  +0.01%  "
